$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 09:48:39"
$wsZh.Range("H2").Value = "2016-03-11 09:48:56"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 09:48:42"
$wsDe.Range("H2").Value = "2016-03-11 09:49:02"
